# Scheduled runner update: refresh market-board derived price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) in columns H:N of the Leve
# tracking tables on each job sheet, per the latest pull of market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2082.3333
$ws.Range("I2").Value = 498
$ws.Range("J2").Value = 2399.2
$ws.Range("K2").Value = 498
$ws.Range("L2").Value = 2399.2
$ws.Range("M2").Value = -385
$ws.Range("N2").Value = -2625.2
$ws.Range("H11").Value = 561.75
$ws.Range("I11").Value = 561.75
$ws.Range("K11").Value = 561.75
$ws.Range("M11").Value = -421.75
$ws.Range("H58").Value = 3779.4443
$ws.Range("J58").Value = 4000
$ws.Range("L58").Value = 12000
$ws.Range("N58").Value = -12300
$ws.Range("H74").Value = 14334.667
$ws.Range("I74").Value = 12000
$ws.Range("K74").Value = 12000
$ws.Range("M74").Value = -11064
$ws.Range("H77").Value = 14334.667
$ws.Range("I77").Value = 12000
$ws.Range("K77").Value = 60000
$ws.Range("M77").Value = -55320

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5839.244
$ws.Range("I32").Value = 4416.817
$ws.Range("K32").Value = 4416.817
$ws.Range("M32").Value = -4129.817
$ws.Range("H61").Value = 2431.1667
$ws.Range("I61").Value = 1965.0182
$ws.Range("J61").Value = 7558.8
$ws.Range("K61").Value = 1965.0182
$ws.Range("L61").Value = 7558.8
$ws.Range("M61").Value = -1753.0182
$ws.Range("N61").Value = -7982.8
$ws.Range("H74").Value = 3034.4
$ws.Range("I74").Value = 2142.3
$ws.Range("J74").Value = 6602.8
$ws.Range("K74").Value = 2142.3
$ws.Range("L74").Value = 6602.8
$ws.Range("M74").Value = -1268.3
$ws.Range("N74").Value = -8350.799999999999
$ws.Range("H77").Value = 3034.4
$ws.Range("I77").Value = 2142.3
$ws.Range("J77").Value = 6602.8
$ws.Range("K77").Value = 10711.5
$ws.Range("L77").Value = 33014
$ws.Range("M77").Value = -6343.5
$ws.Range("N77").Value = -41750
$ws.Range("H86").Value = 44995
$ws.Range("I86").Value = 44995
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 44995
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -43809
$ws.Range("N86").Value = ""
$ws.Range("H89").Value = 44995
$ws.Range("I89").Value = 44995
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 134985
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -129057
$ws.Range("N89").Value = ""
$ws.Range("H136").Value = 2431.1667
$ws.Range("I136").Value = 1965.0182
$ws.Range("J136").Value = 7558.8
$ws.Range("K136").Value = 5895.054599999999
$ws.Range("L136").Value = 22676.4
$ws.Range("M136").Value = -3345.054599999999
$ws.Range("N136").Value = -27776.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7500
$ws.Range("I86").Value = 7500
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 7500
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -6377
$ws.Range("N86").Value = ""
$ws.Range("H88").Value = 750000
$ws.Range("J88").Value = 750000
$ws.Range("L88").Value = 750000
$ws.Range("N88").Value = -750812
$ws.Range("H89").Value = 7500
$ws.Range("I89").Value = 7500
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 37500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -31884
$ws.Range("N89").Value = ""
$ws.Range("H91").Value = 750000
$ws.Range("J91").Value = 750000
$ws.Range("L91").Value = 750000
$ws.Range("N91").Value = -752808
$ws.Range("H94").Value = 1445.1666
$ws.Range("I94").Value = 1484.6538
$ws.Range("K94").Value = 1484.6538
$ws.Range("M94").Value = -1033.6538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2607.2104
$ws.Range("J16").Value = 3496.125
$ws.Range("L16").Value = 3496.125
$ws.Range("N16").Value = -4070.125
$ws.Range("H22").Value = 10162
$ws.Range("I22").Value = 5324
$ws.Range("K22").Value = 5324
$ws.Range("M22").Value = -4974
$ws.Range("H31").Value = 30611.514
$ws.Range("I31").Value = 1558.4
$ws.Range("J31").Value = 48769.707
$ws.Range("K31").Value = 1558.4
$ws.Range("L31").Value = 48769.707
$ws.Range("M31").Value = -1263.4
$ws.Range("N31").Value = -49359.707
$ws.Range("H34").Value = 30611.514
$ws.Range("I34").Value = 1558.4
$ws.Range("J34").Value = 48769.707
$ws.Range("K34").Value = 1558.4
$ws.Range("L34").Value = 48769.707
$ws.Range("M34").Value = -1356.4
$ws.Range("N34").Value = -49173.707
$ws.Range("H58").Value = 3401.9143
$ws.Range("I58").Value = 1901.5714
$ws.Range("K58").Value = 1901.5714
$ws.Range("M58").Value = -1698.5714
$ws.Range("H113").Value = 2607.2104
$ws.Range("J113").Value = 3496.125
$ws.Range("L113").Value = 3496.125
$ws.Range("N113").Value = -7836.125
$ws.Range("H136").Value = 3401.9143
$ws.Range("I136").Value = 1901.5714
$ws.Range("K136").Value = 5704.7142
$ws.Range("M136").Value = -3154.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1963716.8
$ws.Range("I5").Value = 895.63635
$ws.Range("K5").Value = 2686.90905
$ws.Range("M5").Value = -2574.90905
$ws.Range("H10").Value = 48.6
$ws.Range("I10").Value = 29.5
$ws.Range("J10").Value = 125
$ws.Range("K10").Value = 88.5
$ws.Range("L10").Value = 375
$ws.Range("M10").Value = 50.5
$ws.Range("N10").Value = -653
$ws.Range("H131").Value = 7607612.5
$ws.Range("I131").Value = 31250874
$ws.Range("J131").Value = 5118848
$ws.Range("K131").Value = 93752622
$ws.Range("L131").Value = 15356544
$ws.Range("M131").Value = -93747582
$ws.Range("N131").Value = -15366624
$ws.Range("H135").Value = 1963716.8
$ws.Range("I135").Value = 895.63635
$ws.Range("K135").Value = 8060.72715
$ws.Range("M135").Value = -5525.72715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 340193.2
$ws.Range("I80").Value = 626954.9
$ws.Range("K80").Value = 626954.9
$ws.Range("M80").Value = -625956.9
$ws.Range("H83").Value = 340193.2
$ws.Range("I83").Value = 626954.9
$ws.Range("K83").Value = 3134774.5
$ws.Range("M83").Value = -3129782.5
$ws.Range("H102").Value = 2771.5386
$ws.Range("I102").Value = 1213.0526
$ws.Range("K102").Value = 1213.0526
$ws.Range("M102").Value = 408.9474
$ws.Range("H113").Value = 3186.7827
$ws.Range("J113").Value = 4118.5
$ws.Range("L113").Value = 4118.5
$ws.Range("N113").Value = -8458.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2876.158
$ws.Range("I82").Value = 1444.8182
$ws.Range("J82").Value = 4844.25
$ws.Range("K82").Value = 1444.8182
$ws.Range("L82").Value = 4844.25
$ws.Range("M82").Value = -1083.8182
$ws.Range("N82").Value = -5566.25
$ws.Range("H85").Value = 2876.158
$ws.Range("I85").Value = 1444.8182
$ws.Range("J85").Value = 4844.25
$ws.Range("K85").Value = 1444.8182
$ws.Range("L85").Value = 4844.25
$ws.Range("M85").Value = -196.8181999999999
$ws.Range("N85").Value = -7340.25
$ws.Range("H136").Value = 4531.6587
$ws.Range("I136").Value = 3538.1924
$ws.Range("K136").Value = 10614.5772
$ws.Range("M136").Value = -8064.5772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2962.5454
$ws.Range("I132").Value = 2265.111
$ws.Range("K132").Value = 6795.333
$ws.Range("M132").Value = -4265.333
